$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold, border, centered) from the existing last
# header cell (AC1) onto the three new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every data row (2-51): every team in this
# workbook shares the same 1998 season record.
$ws.Range("AD2:AD51").Value = 88
$ws.Range("AE2:AE51").Value = 74
$ws.Range("AF2:AF51").Value = 0

Write-Output "done"
